# Add a mandatory "PAN *" column right after "Name *" (new column B),
# pushing Tags/Category */City one column to the right, and fill in the
# PAN value for each investor row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column A's width so the newly inserted column B can match it
# (mirrors Excel's "insert copied cells" width carry-over behaviour).
$colAWidth = $ws.Columns.Item(1).ColumnWidth

# Insert a new, blank column at B - existing B/C/D (Tags/Category */City)
# shift right to C/D/E.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = $colAWidth

# Header
$ws.Range("B1").Value = "PAN *"

# Data rows - one PAN per investor
$ws.Range("B2").Value = "BUHNXDFEA6"
$ws.Range("B3").Value = "JN2GOV5FYI"
$ws.Range("B4").Value = "CGKT9ROWB1"
$ws.Range("B5").Value = "4I3FNDATK0"
$ws.Range("B6").Value = "5AM81UTOQB"
$ws.Range("B7").Value = "QNEL3S7Z2J"

# Match the editor's final selection
$ws.Range("B7").Select()
